$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("M36").Value = "YES"
$ws.Range("M37").Value = "YES"
$ws.Range("M38").Value = "YES"
$ws.Range("M39").Value = "YES"
$ws.Range("M41").Value = "YES"
$ws.Range("M42").Value = "YES"

$ws.Range("F37:H37").Select()
